# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the "b.md" file
# has now been handed off (zh-cn / de-de), and is no longer a content
# duplicate in zh-cn, while the de-de handback for "b.md" is not yet the
# latest version relative to the current source.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Overview sheet: row 3 (b.md) now shows "Ready for handoff" status for
# both locales, with an updated generation timestamp.
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-17 20:35:49"

# ----------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) status + content-duplicate flag + latest
# handoff file/datetime + error detail. Also widen column P (Error
# Detail) to fit the new long message.
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(16).ColumnWidth = 274 / 7

$zhcn.Range("C3").Value = "Ready for handoff"
# A bare "False" is auto-coerced to a Boolean by Excel's type inference;
# prefix with an apostrophe to force literal text, matching the source
# file's plain shared-string cell, then restore the default cell style
# that the text-coercion path otherwise leaves tagged with a quote-prefix.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-17 20:35:44"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/254d1d5c5964a7ad5183e2a4fa816083943b9e03/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7915e11b22537421692be5862881a4092db5fc78/e2e/b.md."

# ----------------------------------------------------------------------
# de-de sheet: row 3 (b.md) status + latest handoff file/datetime + error
# detail. Widen column P (Error Detail) here too. (Row 2's Latest
# Handback columns keep their existing values - no textual change there.)
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(16).ColumnWidth = 274 / 7

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-17 20:35:49"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/254d1d5c5964a7ad5183e2a4fa816083943b9e03/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7915e11b22537421692be5862881a4092db5fc78/e2e/b.md."
